$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells keep their existing text formatting so that
# numeric-looking strings (e.g. "0.0900", "1.00") are not silently coerced
# into numbers and lose their exact textual representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.885.69"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.351.91"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "544.15"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "136.89"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -5.35%  "
$ws.Range("D9").Value = "2.349.47"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "24.71"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "2.777.32"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "60.855.28"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "2.345.00"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "10.66"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "319.07"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "63.37"
$ws.Range("E25").Value = "  -5.89%  "
$ws.Range("D26").Value = "8.35"
$ws.Range("E26").Value = "  +8.15%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "2.467.84"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "499.02"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "18.46"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +6.64%  "
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "143.35"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "40.61"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "143.12"
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("D46").Value = "3.56"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -8.34%  "
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "19.07"
$ws.Range("E49").Value = "  -6.63%  "
$ws.Range("D50").Value = "0.569"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.0900"
$ws.Range("E51").Value = "  -1.69%  "
